$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the two "Test Datetime" values (shared-string text edits) ---
$ws.Range("E2").Value = "09 Sep 2020 21:50:25"
$ws.Range("E3").Value = "09 Sep 2020 21:52:06"

# --- Re-create the Email hyperlinks on D2/D3 so they get fresh relationship ids ---
# (Existing hyperlink entries loaded from the file can't be edited/removed in place,
#  so clear the whole collection first and re-add - this mirrors the report being
#  regenerated instead of the old hyperlink relationships being silently reused.)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:automationAssessment@ilabquality.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:automationAssessment@ilabquality.com")
